# Apply a row-shuffle to the stats table (columns C:F, rows 2-14) on the
# active worksheet ("Hardik Pandya ").
#
# The new value for each destination row is taken from the source row shown
# in the mapping below (destination row -> source row), derived from the
# unified diff of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original C:F values (rows 2-14) before overwriting anything.
# NOTE: use .Value2 (not .Value) for reading - the .Value getter on this
# runtime is a parameterized property that doesn't resolve correctly.
$original = @{}
for ($r = 2; $r -le 14; $r++) {
    $original[$r] = @(
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2
    )
}

# destination row -> source row
$mapping = @{
    2  = 10
    3  = 6
    4  = 2
    5  = 13
    6  = 8
    7  = 11
    8  = 5
    9  = 14
    10 = 9
    11 = 3
    12 = 7
    13 = 4
    14 = 12
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $vals = $original[$srcRow]
    $ws.Cells.Item($destRow, 3).Value = $vals[0]
    $ws.Cells.Item($destRow, 4).Value = $vals[1]
    $ws.Cells.Item($destRow, 5).Value = $vals[2]
    $ws.Cells.Item($destRow, 6).Value = $vals[3]
}
